$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentor")
$ws.Range("A1:A32").Formula = '=CONCATENATE(C1," ",D1)'
